$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "saya"
$ws.Range("C2").Value = "GBL"
$ws.Range("D2").Value = "2024-10-23 12:15:00"

$ws.Range("E6").Select()
